$wb = $excel.ActiveWorkbook

$newNames = @(
  "summ28910131",
  "summ29103329",
  "summ29322842",
  "summ29534893",
  "summ29750295",
  "summ29969537",
  "summ30188396",
  "summ30406552",
  "summ30619939",
  "summ30839938",
  "summ31061748",
  "summ31363977",
  "summ31782430",
  "summ32251022",
  "summ32558621",
  "summ32863263",
  "summ33177076",
  "summ33473505",
  "summ33780277",
  "summ34083056",
  "summ34375078",
  "summ34671571",
  "summ34957530",
  "summ35254152",
  "summ35540323",
  "summ35812711",
  "summ36091432",
  "summ36390697",
  "summ36673876",
  "summ36954932",
  "summ37238284",
  "summ37540561",
  "summ37825722",
  "summ38107172",
  "summ38390723",
  "summ38674000",
  "summ38957421",
  "summ39238198",
  "summ39525438",
  "summ39805574",
  "summ40090698",
  "summ40373923",
  "summ40654827",
  "summ40940745",
  "summ41223876",
  "summ41527001",
  "summ41821555",
  "summ42095649",
  "summ42373956",
  "summ42657277"
)

for ($i = 0; $i -lt $newNames.Length; $i++) {
    $idx = $i + 1
    $ws = $wb.Worksheets.Item($idx)
    $ws.Name = $newNames[$i]
}
